$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "RM 232" (original row 26) and "SC 92" (original row 28) were removed from the
# dataset, shifting all subsequent rows up. Delete them first (26, then the new 27,
# which is the old row 28 after the first shift).
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-point the missing-data mask: clear cells that are now considered missing and
# restore/fix the underlying numeric values for cells that are no longer missing.
$ws.Range("D2").ClearContents()

$ws.Range("F4").ClearContents()

$ws.Range("D5").Value = -14.4

$ws.Range("C6").Value = 15.1
$ws.Range("D6").Value = -14.2
$ws.Range("F6").Value = 16.43

$ws.Range("C8").ClearContents()

$ws.Range("D10").ClearContents()
$ws.Range("F10").Value = 16.43

$ws.Range("F11").Value = 17.65

$ws.Range("C12").Value = 12.5
$ws.Range("F12").ClearContents()

$ws.Range("D13").ClearContents()

$ws.Range("C14").ClearContents()
$ws.Range("F14").Value = 17.76

$ws.Range("F16").ClearContents()

$ws.Range("C17").Value = 11.2
$ws.Range("F17").ClearContents()

$ws.Range("C18").Value = 11.5

$ws.Range("C19").ClearContents()
$ws.Range("F19").Value = 17.81

$ws.Range("C20").ClearContents()

$ws.Range("F21").Value = 16.58

$ws.Range("F22").Value = 16.81

$ws.Range("C23").Value = 12.2

$ws.Range("D24").Value = -13.9

$ws.Range("F25").ClearContents()

$ws.Range("F26").ClearContents()

$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()
$ws.Range("F27").ClearContents()

$ws.Range("D28").ClearContents()
$ws.Range("F28").ClearContents()

$ws.Range("B29").ClearContents()

$ws.Range("D30").Value = -13.6

$ws.Range("F31").Value = 17.18

$ws.Range("B32").ClearContents()
